$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A and B by 1 character (15.42578125 -> 16.42578125).
# Note: the COM layer here quantizes ColumnWidth to the nearest 1/6 of a
# character, so 15.666666666666668 is the closest settable value that
# rounds (on save) to the desired stored width of ~16.43.
$ws.Columns.Item(1).ColumnWidth = 15.666666666666668
$ws.Columns.Item(2).ColumnWidth = 15.666666666666668

$values = @(
    @(-0.2476608498783861, 0.24726123075171813),
    @(-0.20115084388407745, 0.19972567715326761),
    @(-0.096781991906656373, 0.096551037798221273),
    @(-0.15854683200707598, 0.15792952456073195),
    @(-0.15192952495581835, 0.15070219646036875),
    @(-0.079135062929207134, 0.079058005208028526),
    @(-0.059058005691547066, 0.058885182443944473),
    @(-0.038885182931523765, 0.03873899006189685),
    @(-0.032738990478989649, 0.032613100911402348),
    @(-0.026613101331804501, 0.026596436273152335),
    @(-0.037023829609406533, 0.036989013698242701),
    @(-0.030989014119871872, 0.030890584517594188),
    @(-0.024890584943673133, 0.024868008868227598),
    @(-0.012868009327418051, 0.012861185065538017),
    @(-0.0068611854934941263, 0.006857728622906123),
    @(-0.00085772905127523558, 0.00085255451655452674),
    @(-0.0090043415338634247, 0.0089999995551810485),
    @(-0.036110387266482746, 0.03609674842894961),
    @(-0.027096748822755146, 0.0270136953078155),
    @(-0.018013695705388244, 0.018004289101750004),
    @(-0.0090042894998889622, 0.008999999601473796),
    @(-0.093943685649238162, 0.093631934958597896),
    @(-0.084631935364845035, 0.084126307112511789),
    @(-0.042126307701066068, 0.041999999408152355),
    @(-0.065567133938898792, 0.065510013455547522),
    @(-0.011327512257359729, 0.011282167807969756),
    @(-0.0052821682164667649, 0.0051287225838825279),
    @(0.00087127700711775447, -0.00097652320257957115),
    @(0.012976522761793063, -0.013028199027992571),
    @(-0.042163060215933523, 0.042019076826638635),
    @(-0.027019077286691839, 0.027000792563438125),
    @(-0.0060007930560495026, 0.0059999995870079204)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
